$wb = $excel.ActiveWorkbook

# --- Update statistics sheets: one more student became "Reprobado" instead of "Blanco" ---
$ws1P = $wb.Worksheets.Item("Estadisticos 1P")
$ws1P.Range("D2").Value = 0
$ws1P.Range("E2").Value = 10

$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 17

$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("D2").Value = 0
$wsFinal.Range("E2").Value = 17

# --- Add two new students to the "Rescatables" sheet ---
$wsResc = $wb.Worksheets.Item("Rescatables")

$wsResc.Range("A2").Value = 20330051920178
$wsResc.Range("B2").Value = "QUIRIZ"
$wsResc.Range("C2").Value = "RAMOS"
$wsResc.Range("D2").Value = "MONICA"
$wsResc.Range("E2").Value = "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA"
$wsResc.Range("F2").Value = "6ARHM"
$wsResc.Range("G2").Value = 2

$wsResc.Range("A3").Value = 20330051920373
$wsResc.Range("B3").Value = "RICO"
$wsResc.Range("C3").Value = "BAUTISTA"
$wsResc.Range("D3").Value = "EDGAR RAMSES"
$wsResc.Range("E3").Value = "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA"
$wsResc.Range("F3").Value = "6ARHM"
$wsResc.Range("G3").Value = 2
